# Update "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets
# to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览"     = @{ "F4" = 7788; "F8" = 6489; "F13" = 36; "F15" = 56; "F21" = 3773; "F30" = 1731 }
    "全部类型" = @{ "F7" = 7788; "F11" = 6489; "F14" = 36; "F15" = 56; "F22" = 3773; "F33" = 1731 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $updates[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = $cellMap[$cellRef]
    }
}
